$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-of dates) ---
$ws.Range("A8").Value = "Volume 32   Number  21"
$ws.Range("C9").Value = "Report Covering the Week  5/19/2025  Through  5/25/2025"

# --- Row 15 ---
$ws.Range("F15").Value = 3
$ws.Range("H15").Value = 50

# --- Row 16 ---
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 10
$ws.Range("E16").Value = -40
$ws.Range("G16").Value = 28
$ws.Range("H16").Value = -7.142857142857
$ws.Range("I16").Value = 116
$ws.Range("J16").Value = 145
$ws.Range("K16").Value = -20
$ws.Range("L16").Value = 10.476190476190
$ws.Range("M16").Value = -15.328467153284
$ws.Range("N16").Value = -77.029702970297

# --- Row 17 ---
$ws.Range("D17").Value = 17
$ws.Range("E17").Value = -47.058823529411
$ws.Range("F17").Value = 38
$ws.Range("G17").Value = 51
$ws.Range("H17").Value = -25.490196078431
$ws.Range("I17").Value = 209
$ws.Range("J17").Value = 218
$ws.Range("K17").Value = -4.128440366972
$ws.Range("L17").Value = 10
$ws.Range("M17").Value = 74.166666666666

# --- Row 18 ---
$ws.Range("D18").Value = 4
$ws.Range("E18").Value = -50
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 15
$ws.Range("H18").Value = -60
$ws.Range("I18").Value = 62
$ws.Range("J18").Value = 87
$ws.Range("K18").Value = -28.735632183908
$ws.Range("L18").Value = 8.771929824561
$ws.Range("M18").Value = -53.030303030303
$ws.Range("N18").Value = -92.757009345794

# --- Row 19 ---
$ws.Range("C19").Value = 12
$ws.Range("D19").Value = 20
$ws.Range("E19").Value = -40
$ws.Range("F19").Value = 43
$ws.Range("G19").Value = 62
$ws.Range("H19").Value = -30.645161290322
$ws.Range("I19").Value = 269
$ws.Range("J19").Value = 355
$ws.Range("K19").Value = -24.225352112676
$ws.Range("L19").Value = -13.504823151125
$ws.Range("M19").Value = 55.491329479768
$ws.Range("N19").Value = -48.269230769230

# --- Row 20 ---
$ws.Range("C20").Value = 9
$ws.Range("E20").Value = 28.571428571428
$ws.Range("F20").Value = 31
$ws.Range("G20").Value = 29
$ws.Range("H20").Value = 6.896551724137
$ws.Range("I20").Value = 100
$ws.Range("J20").Value = 117
$ws.Range("K20").Value = -14.529914529914
$ws.Range("L20").Value = -20
$ws.Range("M20").Value = 5.263157894736
$ws.Range("N20").Value = -88.425925925925

# --- Row 21 ---
$ws.Range("C21").Value = 38
$ws.Range("D21").Value = 58
$ws.Range("E21").Value = -34.482758620689
$ws.Range("F21").Value = 148
$ws.Range("G21").Value = 187
$ws.Range("H21").Value = -20.855614973262
$ws.Range("I21").Value = 774
$ws.Range("J21").Value = 934
$ws.Range("K21").Value = -17.130620985010
$ws.Range("L21").Value = -2.641509433962
$ws.Range("M21").Value = 16.041979010494
$ws.Range("N21").Value = -73.447684391080

# --- Row 22 ---
$ws.Range("D22").Value = 4
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 10
$ws.Range("H22").Value = -60
$ws.Range("I22").Value = 24
$ws.Range("J22").Value = 31
$ws.Range("K22").Value = -22.580645161290
$ws.Range("L22").Value = -47.826086956521
$ws.Range("M22").Value = 71.428571428571

# --- Row 24 ---
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 40
$ws.Range("E24").Value = -55
$ws.Range("F24").Value = 92
$ws.Range("G24").Value = 145
$ws.Range("H24").Value = -36.551724137931
$ws.Range("I24").Value = 516
$ws.Range("J24").Value = 926
$ws.Range("K24").Value = -44.276457883369
$ws.Range("L24").Value = -37.530266343825
$ws.Range("M24").Value = 24.038461538461

# --- Row 25 ---
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 21
$ws.Range("E25").Value = -66.666666666666
$ws.Range("F25").Value = 29
$ws.Range("G25").Value = 70
$ws.Range("H25").Value = -58.571428571428
$ws.Range("I25").Value = 214
$ws.Range("J25").Value = 522
$ws.Range("K25").Value = -59.003831417624
$ws.Range("L25").Value = -51.363636363636

# --- Row 26 ---
$ws.Range("C26").Value = 24
$ws.Range("D26").Value = 26
$ws.Range("E26").Value = -7.692307692307
$ws.Range("F26").Value = 94
$ws.Range("G26").Value = 116
$ws.Range("H26").Value = -18.965517241379
$ws.Range("I26").Value = 437
$ws.Range("J26").Value = 494
$ws.Range("K26").Value = -11.538461538461
$ws.Range("L26").Value = 24.857142857142
$ws.Range("M26").Value = 19.398907103825

# --- Row 27 ---
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = 0

# --- Row 28 ---
$ws.Range("F28").Value = 19
$ws.Range("G28").Value = 8
$ws.Range("H28").Value = 137.5
$ws.Range("I28").Value = 51
$ws.Range("J28").Value = 56
$ws.Range("K28").Value = -8.928571428571
$ws.Range("L28").Value = -8.928571428571

# --- Row 31 --- (F31 becomes numeric, was a text "0" placeholder style)
$ws.Range("F31").NumberFormat = "#,##0"
$ws.Range("F31").Value = 1
$ws.Range("I31").Value = 3
$ws.Range("K31").Value = -40
$ws.Range("L31").Value = -25
